$d = $word.ActiveDocument

# First paragraph: the one holding the "**ID__AFFARS_...__ID**" placeholder run
$p1 = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right) with 5pt spacing from text,
# matching the new <w:pBdr> with w:space="5" on every edge.
$p1.Format.Borders.DistanceFromTop = 5
$p1.Format.Borders.DistanceFromLeft = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromRight = 5

# Update the left indent from 120 twips to 225 twips (LeftIndent is in points).
$p1.Format.LeftIndent = 225 / 20

# Replace the placeholder text and drop the trailing lone-space run in one go:
# searching across both runs' text ("...__ID**" + " ") and replacing with the
# new id text (no trailing space) merges/removes the second run.
$null = $d.Content.Find.Execute("**ID__AFFARS_pgi_5306_topic_4__ID** ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "**ID__AFFARS_AFMC_PGI_5306__ID**", 2)
